# Commit: "Fruta / hortaliza, semanal"
#
# The sheet gains one new daily record. A new row is inserted at row 30
# (pushing the former rows 30-92 down to 31-93, so the old last row, 92,
# becomes the new row 93), and the freshly inserted row 30 is populated
# with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 30; this shifts rows 30..92
# down to 31..93 and grows the sheet dimension to A1:R93, matching the
# diff exactly.
$ws.Rows.Item(30).Insert()

# Populate the newly-inserted (now blank) row 30 with the new record.
$ws.Range("A30").Value() = 10
$ws.Range("B30").Value() = "Vega Modelo de Temuco"
$ws.Range("C30").Value() = "La Araucanía"
$ws.Range("D30").Value() = 44498
$ws.Range("E30").Value() = 9
$ws.Range("F30").Value() = 100112012
$ws.Range("G30").Value() = "Espinaca"
$ws.Range("H30").Value() = "Sin especificar"
$ws.Range("I30").Value() = "Primera"
$ws.Range("J30").Value() = 40
$ws.Range("K30").Value() = 7000
$ws.Range("L30").Value() = 7000
$ws.Range("M30").Value() = 7000
$ws.Range("N30").Value() = "$/docena de atados"
$ws.Range("O30").Value() = "Región de La Araucanía"
$ws.Range("P30").Value() = 2333
$ws.Range("Q30").Value() = 3
$ws.Range("R30").Value() = "Hortaliza"
